$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "new dataset / new model" values for A2:A49
$values = @(
    64.769647696476966,
    64.769647696476966,
    65.040650406504056,
    64.498644986449861,
    63.956639566395665,
    64.498644986449861,
    65.853658536585371,
    68.834688346883468,
    69.647696476964768,
    67.208672086720867,
    65.311653116531161,
    68.563685636856363,
    69.105691056910572,
    68.834688346883468,
    68.834688346883468,
    67.208672086720867,
    67.750677506775077,
    71.002710027100264,
    66.937669376693762,
    66.124661246612476,
    65.853658536585371,
    59.078590785907856,
    59.078590785907856,
    57.72357723577236,
    67.750677506775077,
    67.208672086720867,
    69.918699186991873,
    68.563685636856363,
    68.021680216802167,
    69.105691056910572,
    59.078590785907856,
    59.891598915989164,
    60.433604336043359,
    57.452574525745263,
    55.826558265582662,
    60.975609756097562,
    56.36856368563685,
    56.36856368563685,
    56.36856368563685,
    67.750677506775077,
    68.834688346883468,
    69.376693766937663,
    66.937669376693762,
    65.311653116531161,
    67.208672086720867,
    65.040650406504056,
    66.395663956639567,
    67.750677506775077
)

# Header label (unchanged text, but the source workbook's shared-string
# table gained a duplicate entry for it when it was regenerated for the
# new model/dataset run).
$ws.Range("A1").Value = "HK_G_acc_SD"

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
